$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# The source file stores every cell (including numeric-looking ones like
# "10") as text. Writing a plain numeric-looking string via .Value gets
# auto-coerced to a real number by this engine, so force each target cell
# to stay text by switching it to a text number format before the write,
# then restoring the Normal style so it isn't visually different from the
# rest of the sheet (applied per-cell - multi-area Range unions only
# re-style their first area here, so each cell is handled individually).
function Set-TextCell($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# New order rows appended after the existing data (rows 2-61), extending
# the sheet from A1:L61 to A1:L71.
Set-TextCell "A62" "1"
Set-TextCell "C62" "478_绿芯向日葵_sunflower mini_undefined_1bunch"
Set-TextCell "F62" "10"

Set-TextCell "C63" "521_商陆_phytolacca acinosa _undefined_1bunch"
Set-TextCell "F63" "15"

Set-TextCell "C64" "328_卢荀草_undefined_undefined_1bunch"
Set-TextCell "F64" "10"

Set-TextCell "C65" "203_佛罗伊德_Floyd_Rosa rugosa Thunb._20stems"
Set-TextCell "F65" "16"

Set-TextCell "C66" "277_草莓杏仁饼_undefined_Rosa rugosa Thunb._10stems"
Set-TextCell "F66" "6"

Set-TextCell "C67" "224_折射_Reflex_Rosa rugosa Thunb._10stems"
Set-TextCell "F67" "5"

Set-TextCell "C68" "404_小飞燕白色_ delphinium ballkleid`nwhite_undefined_1bunch"
Set-TextCell "F68" "10"

Set-TextCell "C69" "405_小飞燕浅蓝_ delphinium ballkleid`ndark blue_undefined_1bunch"
Set-TextCell "F69" "10"

Set-TextCell "C70" "497_小飞燕粉色_delphinium ballkleid`npink_undefined_1bunch"
Set-TextCell "F70" "10"

Set-TextCell "A71" "2"

# Summary sheet: the rolled-up per-package counts string grows to reflect
# the newly appended orders.
$summary = $wb.Worksheets.Item("Summary")
$g2 = $summary.Range("G2")
$g2.NumberFormat = "@"
$g2.Value = "0151540401033532151014713101491410105510115111082615151515151041595010201555510101551055151051010101010151016651010100"
$g2.Style = "Normal"
